$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and reorder rows 47-48) per latest scrape

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.256.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.37%  "

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.689.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.14%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5264"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.46%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.05%  "

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2699"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.08%  "

# Row 9: Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06448"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.63%  "

# Row 10: Solana
$ws.Range("E10").Value = "  +2.42%  "

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07470"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.60%  "

# Row 12: WrappedEther
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.696.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.49%  "

# Row 13: Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.556"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.02%  "

# Row 14: Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5864"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.38%  "

# Row 15: ShibaInu
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008556"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.57%  "

# Row 16: Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.60%  "

# Row 17: WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.305.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.28%  "

# Row 18: Uniswap
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.971"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.15%  "

# Row 19: Dai
$ws.Range("E19").Value = "  +0.09%  "

# Row 20: Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.32%  "

# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.88%  "

# Row 22: Chainlink
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.245"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.69%  "

# Row 23: BinanceUSD
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24: Monero
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.29%  "

# Row 25: Cosmos
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.689"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "

# Row 26: Stellar
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1237"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.82%  "

# Row 27: EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.82%  "

# Row 28: Hedera
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06709"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.58%  "

# Row 29: Toncoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.356"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.08%  "

# Row 30: PancakeSwap
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.330"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.76%  "

# Row 31: Filecoin
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.597"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.76%  "

# Row 32: InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.557"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.36%  "

# Row 33: LidoDAOToken
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.671"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.69%  "

# Row 34: ARBITRUM
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.029"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.05%  "

# Row 35: ImmutableX
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6217"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.48%  "

# Row 36: HuobiToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.385"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.29%  "

# Row 37: MXToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.710"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.48%  "

# Row 38: FraxShare
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.289"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.93%  "

# Row 39: VeChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01620"
$ws.Range("D39").Style = "Normal"

# Row 40: Maker
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.103.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.04%  "

# Row 41: TrustWalletToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8874"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.46%  "

# Row 42: PaxDollar
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.015"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.77%  "

# Row 43: Quant
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.68%  "

# Row 44: RocketPoolETH
$ws.Range("E44").Value = "  +0.89%  "

# Row 45: BabyDogeCoin
$ws.Range("E45").Value = "  +3.00%  "

# Row 46: Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.84%  "

# Row 47: EnergySwap
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.176"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.16%  "

# Row 48: Frax
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "

# Row 49: Cronos
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05261"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.55%  "

# Row 50: Mantle
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4294"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.07%  "

# Row 51: Aptos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.043"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.02%  "

